$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '68.748.70'
Set-TextValue $ws.Range('E2') '  +1.41%  '
Set-TextValue $ws.Range('D3') '3.751.69'
Set-TextValue $ws.Range('E3') '  -1.35%  '
Set-TextValue $ws.Range('E4') '  +0.06%  '
Set-TextValue $ws.Range('D5') '602.29'
Set-TextValue $ws.Range('E5') '  -0.04%  '
Set-TextValue $ws.Range('D6') '169.22'
Set-TextValue $ws.Range('E6') '  -1.14%  '
Set-TextValue $ws.Range('D7') '3.748.62'
Set-TextValue $ws.Range('E7') '  -1.36%  '
Set-TextValue $ws.Range('E8') '  -0.07%  '
Set-TextValue $ws.Range('E9') '  +0.98%  '
Set-TextValue $ws.Range('D10') '0.164'
Set-TextValue $ws.Range('E10') '  +2.97%  '
Set-TextValue $ws.Range('D11') '6.34'
Set-TextValue $ws.Range('E11') '  +2.33%  '
Set-TextValue $ws.Range('D12') '0.462'
Set-TextValue $ws.Range('E12') '  -1.21%  '
Set-TextValue $ws.Range('D13') '38.27'
Set-TextValue $ws.Range('E13') '  -1.55%  '
Set-TextValue $ws.Range('E14') '  +0.59%  '
Set-TextValue $ws.Range('D15') '4.384.62'
Set-TextValue $ws.Range('E15') '  -1.10%  '
Set-TextValue $ws.Range('D16') '3.754.30'
Set-TextValue $ws.Range('E16') '  -1.29%  '
Set-TextValue $ws.Range('D17') '68.766.00'
Set-TextValue $ws.Range('D18') '7.30'
Set-TextValue $ws.Range('E18') '  +0.48%  '
Set-TextValue $ws.Range('E19') '  -0.17%  '
Set-TextValue $ws.Range('D20') '17.17'
Set-TextValue $ws.Range('E20') '  -1.72%  '
Set-TextValue $ws.Range('D21') '10.72'
Set-TextValue $ws.Range('E21') '  +12.67%  '
Set-TextValue $ws.Range('D22') '495.48'
Set-TextValue $ws.Range('E22') '  +0.04%  '
Set-TextValue $ws.Range('D23') '0.729'
Set-TextValue $ws.Range('E23') '  -2.22%  '
Set-TextValue $ws.Range('D24') '85.55'
Set-TextValue $ws.Range('E24') '  -0.53%  '
Set-TextValue $ws.Range('E25') '  +1.11%  '
Set-TextValue $ws.Range('D26') '2.31'
Set-TextValue $ws.Range('E26') '  -3.52%  '
Set-TextValue $ws.Range('D27') '12.41'
Set-TextValue $ws.Range('E27') '  +0.28%  '
Set-TextValue $ws.Range('D28') '10.17'
Set-TextValue $ws.Range('E28') '  -0.57%  '
Set-TextValue $ws.Range('E29') '  -0.22%  '
Set-TextValue $ws.Range('E30') '  +3.68%  '
Set-TextValue $ws.Range('E31') '  +0.17%  '
Set-TextValue $ws.Range('D32') '7.94'
Set-TextValue $ws.Range('E32') '  +0.60%  '
Set-TextValue $ws.Range('D33') '32.07'
Set-TextValue $ws.Range('E33') '  -3.05%  '
Set-TextValue $ws.Range('D34') '3.898.16'
Set-TextValue $ws.Range('E34') '  -1.22%  '
Set-TextValue $ws.Range('D35') '3.686.96'
Set-TextValue $ws.Range('E35') '  -1.41%  '
Set-TextValue $ws.Range('E36') '  -1.42%  '
Set-TextValue $ws.Range('D37') '1.00'
Set-TextValue $ws.Range('E37') '  +0.08%  '
Set-TextValue $ws.Range('D38') '1.01'
Set-TextValue $ws.Range('E38') '  -0.58%  '
Set-TextValue $ws.Range('E39') '  +0.04%  '
Set-TextValue $ws.Range('E40') '  -0.01%  '
Set-TextValue $ws.Range('D41') '0.326'
Set-TextValue $ws.Range('E41') '  -1.34%  '
Set-TextValue $ws.Range('D42') '437.94'
Set-TextValue $ws.Range('E42') '  -4.76%  '
Set-TextValue $ws.Range('D43') '49.02'
Set-TextValue $ws.Range('E43') '  -0.29%  '
Set-TextValue $ws.Range('D44') '1.98'
Set-TextValue $ws.Range('E44') '  -1.77%  '
Set-TextValue $ws.Range('D45') '2.86'
Set-TextValue $ws.Range('E45') '  -0.14%  '
Set-TextValue $ws.Range('D46') '8.51'
Set-TextValue $ws.Range('E46') '  +0.49%  '
Set-TextValue $ws.Range('E47') '  -0.05%  '
Set-TextValue $ws.Range('D48') '40.70'
Set-TextValue $ws.Range('E48') '  -0.04%  '
Set-TextValue $ws.Range('D49') '2.824.74'
Set-TextValue $ws.Range('E49') '  -0.79%  '
Set-TextValue $ws.Range('D50') '141.04'
Set-TextValue $ws.Range('E50') '  +0.51%  '
Set-TextValue $ws.Range('D51') '0.0355'
Set-TextValue $ws.Range('E51') '  +0.54%  '
